$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column cells keep their original Text format so that
# numeric-looking strings (e.g. "0.994", "60.20") are not coerced to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.204.56"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "2.611.74"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "521.43"
$ws.Range("E5").Value = "  +1.10%  "
$ws.Range("D6").Value = "148.88"
$ws.Range("E6").Value = "  -3.77%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  -5.27%  "
$ws.Range("D9").Value = "2.617.71"
$ws.Range("E9").Value = "  +0.24%  "
$ws.Range("D10").Value = "6.31"
$ws.Range("E10").Value = "  -4.97%  "
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("E12").Value = "  -1.68%  "
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("D14").Value = "3.068.13"
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("D15").Value = "60.221.21"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("E16").Value = "  -2.35%  "
$ws.Range("E17").Value = "  -1.68%  "
$ws.Range("D18").Value = "2.614.05"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("E19").Value = "  -2.78%  "
$ws.Range("D20").Value = "342.94"
$ws.Range("E20").Value = "  -3.04%  "
$ws.Range("D21").Value = "10.41"
$ws.Range("E21").Value = "  -1.93%  "
$ws.Range("E22").Value = "  -2.20%  "
$ws.Range("D23").Value = "0.994"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").Value = "60.20"
$ws.Range("E24").Value = "  -1.40%  "
$ws.Range("E25").Value = "  -2.67%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "0.163"
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "0.996"
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").Value = "0.0₃0805"
$ws.Range("E28").Value = "  -4.47%  "
$ws.Range("D29").Value = "7.07"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").Value = "6.03"
$ws.Range("E31").Value = "  -2.68%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").Value = "18.93"
$ws.Range("E33").Value = "  -2.71%  "
$ws.Range("D34").Value = "149.80"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").Value = "3.96"
$ws.Range("E35").Value = "  -2.98%  "
$ws.Range("D36").Value = "0.916"
$ws.Range("E36").Value = "  -3.10%  "
$ws.Range("E37").Value = "  -4.95%  "
$ws.Range("D38").Value = "0.864"
$ws.Range("E38").Value = "  +2.56%  "
$ws.Range("D39").Value = "36.36"
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("E40").Value = "  -3.20%  "
$ws.Range("E41").Value = "  -4.06%  "
$ws.Range("D42").Value = "287.05"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("E44").Value = "  -1.17%  "
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("E46").Value = "  -2.30%  "
$ws.Range("D47").Value = "19.53"
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("E49").Value = "  +0.95%  "
$ws.Range("E50").Value = "  -5.67%  "
$ws.Range("D51").Value = "1.950.06"
$ws.Range("E51").Value = "  -1.48%  "
